$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 114
$ws.Range("D2").Value = 106
$ws.Range("J2").Value = 21649

$ws.Range("C3").Value = 346
$ws.Range("D3").Value = 456
$ws.Range("J3").Value = 113030

$ws.Range("C4").Value = 358
$ws.Range("D4").Value = 462
$ws.Range("J4").Value = 288016

$ws.Range("C5").Value = 232
$ws.Range("D5").Value = 307
$ws.Range("J5").Value = 144417
